$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2328.5715
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 2466.6667
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 2466.6667
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -3118.6667

$ws.Range("H55").Value = 329.69232
$ws.Range("I55").Value = 287
$ws.Range("J55").Value = 425.75
$ws.Range("K55").Value = 287
$ws.Range("L55").Value = 425.75
$ws.Range("M55").Value = -73
$ws.Range("N55").Value = -853.75

$ws.Range("H100").Value = 250000260
$ws.Range("I100").Value = 250000260
$ws.Range("K100").Value = 250000260
$ws.Range("M100").Value = -249999719

$ws.Range("H107").Value = 2838.3635
$ws.Range("I107").Value = 3502
$ws.Range("J107").Value = 2285.3333
$ws.Range("K107").Value = 3502
$ws.Range("L107").Value = 2285.3333
$ws.Range("M107").Value = -1582
$ws.Range("N107").Value = -6125.3333

$ws.Range("H137").Value = 1131.0769
$ws.Range("I137").Value = 1026.4348
$ws.Range("J137").Value = 1933.3334
$ws.Range("K137").Value = 3079.3044
$ws.Range("L137").Value = 5800.0002
$ws.Range("M137").Value = -529.3044
$ws.Range("N137").Value = -10900.0002

$ws.Range("H141").Value = 1162.5555
$ws.Range("I141").Value = 1044.75
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 3134.25
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 2045.75
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3436.7017
$ws.Range("I32").Value = 3250.549
$ws.Range("K32").Value = 3250.549
$ws.Range("M32").Value = -2963.549

$ws.Range("H45").Value = 1414
$ws.Range("I45").Value = 1508.909
$ws.Range("J45").Value = 892
$ws.Range("K45").Value = 1508.909
$ws.Range("L45").Value = 892
$ws.Range("M45").Value = -1131.909
$ws.Range("N45").Value = -1646

$ws.Range("H61").Value = 1237.7407
$ws.Range("I61").Value = 1072.3636
$ws.Range("J61").Value = 1965.4
$ws.Range("K61").Value = 1072.3636
$ws.Range("L61").Value = 1965.4
$ws.Range("M61").Value = -860.3635999999999
$ws.Range("N61").Value = -2389.4

$ws.Range("H132").Value = 2206.16
$ws.Range("I132").Value = 1817.1
$ws.Range("K132").Value = 5451.299999999999
$ws.Range("M132").Value = -2921.299999999999

$ws.Range("H136").Value = 1237.7407
$ws.Range("I136").Value = 1072.3636
$ws.Range("J136").Value = 1965.4
$ws.Range("K136").Value = 3217.0908
$ws.Range("L136").Value = 5896.200000000001
$ws.Range("M136").Value = -667.0907999999999
$ws.Range("N136").Value = -10996.2

$ws.Range("H139").Value = 32730.666
$ws.Range("J139").Value = 32730.666
$ws.Range("L139").Value = 32730.666
$ws.Range("N139").Value = -43010.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 399.7
$ws.Range("I80").Value = 307.2857
$ws.Range("J80").Value = 449.46155
$ws.Range("K80").Value = 307.2857
$ws.Range("L80").Value = 449.46155
$ws.Range("M80").Value = 690.7143
$ws.Range("N80").Value = -2445.46155

$ws.Range("H83").Value = 399.7
$ws.Range("I83").Value = 307.2857
$ws.Range("J83").Value = 449.46155
$ws.Range("K83").Value = 1536.4285
$ws.Range("L83").Value = 2247.30775
$ws.Range("M83").Value = 3455.5715
$ws.Range("N83").Value = -12231.30775

$ws.Range("H99").Value = 38462984
$ws.Range("I99").Value = 45455956
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 45455956
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -45454458
$ws.Range("N99").Value = -4646

$ws.Range("H105").Value = 250002500
$ws.Range("I105").Value = 250002500
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 250002500
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -250000753
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 6619
$ws.Range("I134").Value = 1089.0667
$ws.Range("J134").Value = 14913.9
$ws.Range("K134").Value = 3267.2001
$ws.Range("L134").Value = 44741.7
$ws.Range("M134").Value = -732.2001
$ws.Range("N134").Value = -49811.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2000.5834
$ws.Range("J69").Value = 2018.8182
$ws.Range("L69").Value = 6056.4546
$ws.Range("N69").Value = -7678.4546

$ws.Range("H72").Value = 2000.5834
$ws.Range("J72").Value = 2018.8182
$ws.Range("L72").Value = 18169.3638
$ws.Range("N72").Value = -26281.3638

$ws.Range("H131").Value = 15385840
$ws.Range("J131").Value = 1327.6492
$ws.Range("L131").Value = 3982.9476
$ws.Range("N131").Value = -14062.9476

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1456.091
$ws.Range("I102").Value = 1502.1875
$ws.Range("J102").Value = 1333.1666
$ws.Range("K102").Value = 1502.1875
$ws.Range("L102").Value = 1333.1666
$ws.Range("M102").Value = 119.8125
$ws.Range("N102").Value = -4577.1666

$ws.Range("H132").Value = 3307.0908
$ws.Range("I132").Value = 2911.4285
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 8734.2855
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -6204.2855
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H132").Value = 20904.5
$ws.Range("I132").Value = 1417.5834
$ws.Range("J132").Value = 37607.57
$ws.Range("K132").Value = 4252.7502
$ws.Range("L132").Value = 112822.71
$ws.Range("M132").Value = -1722.7502
$ws.Range("N132").Value = -117882.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("K88").Value = 10000
$ws.Range("M88").Value = -9594

$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("K91").Value = 10000
$ws.Range("M91").Value = -8596

$ws.Range("H96").Value = 926.7273
$ws.Range("I96").Value = 721.4
$ws.Range("K96").Value = 721.4
$ws.Range("M96").Value = 651.6

$ws.Range("H125").Value = 99990
$ws.Range("J125").Value = 99990
$ws.Range("L125").Value = 99990
$ws.Range("N125").Value = -109830

$ws.Range("H126").Value = 142858220
$ws.Range("I126").Value = 142858220
$ws.Range("K126").Value = 428574660
$ws.Range("M126").Value = -428572190
